# Updates cryptocurrency price/volume data (and reorders rows 44-46)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.095.60"
$ws.Range("E2").Value = "  +0.21%  "

# Row 3
$ws.Range("D3").Value = "1.834.32"
$ws.Range("E3").Value = "  +0.00%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.42%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.47"
$ws.Range("E5").Value = "  +0.46%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6276"
$ws.Range("E6").Value = "  -0.41%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  +0.25%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07508"
$ws.Range("E8").Value = "  -1.37%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2921"
$ws.Range("E9").Value = "  -0.39%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.20"
$ws.Range("E10").Value = "  +2.63%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07676"
$ws.Range("E11").Value = "  -0.54%  "

# Row 12
$ws.Range("D12").Value = "1.838.46"
$ws.Range("E12").Value = "  +0.19%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.005"
$ws.Range("E13").Value = "  +0.88%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6670"
$ws.Range("E14").Value = "  +0.09%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.68"
$ws.Range("E15").Value = "  -0.20%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009374"
$ws.Range("E16").Value = "  -7.63%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.981"
$ws.Range("E17").Value = "  -1.18%  "

# Row 18
$ws.Range("D18").Value = "29.101.72"
$ws.Range("E18").Value = "  +0.19%  "

# Row 19
$ws.Range("D19").Value = "2.080.58"
$ws.Range("E19").Value = "  -0.18%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.58"
$ws.Range("E20").Value = "  +1.73%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "223.15"
$ws.Range("E21").Value = "  -1.68%  "

# Row 22
$ws.Range("E22").Value = "  +0.53%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.097"
$ws.Range("E23").Value = "  -1.26%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.003"
$ws.Range("E24").Value = "  +0.41%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.60"
$ws.Range("E25").Value = "  +0.80%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1397"
$ws.Range("E26").Value = "  +1.39%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.487"
$ws.Range("E27").Value = "  -0.23%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.90"
$ws.Range("E28").Value = "  -0.15%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.498"
$ws.Range("E29").Value = "  +0.73%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05665"
$ws.Range("E30").Value = "  +8.22%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.153"
$ws.Range("E31").Value = "  +0.90%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.075"
$ws.Range("E32").Value = "  +1.37%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.205"
$ws.Range("E33").Value = "  +1.18%  "

# Row 34
$ws.Range("E34").Value = "  -0.29%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7427"
$ws.Range("E35").Value = "  +0.78%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.139"
$ws.Range("E36").Value = "  -0.10%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.670"
$ws.Range("E37").Value = "  -1.21%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.766"
$ws.Range("E38").Value = "  +0.28%  "

# Row 39
$ws.Range("D39").Value = "1.221.67"
$ws.Range("E39").Value = "  -1.51%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01780"
$ws.Range("E40").Value = "  -0.42%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.532"
$ws.Range("E41").Value = "  +2.76%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8940"
$ws.Range("E42").Value = "  -0.05%  "

# Row 43
$ws.Range("E43").Value = "  +0.13%  "

# Row 44
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.91"
$ws.Range("E44").Value = "  +0.16%  "

# Row 45
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.979.03"
$ws.Range("E45").Value = "  -0.17%  "

# Row 46
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000125"
$ws.Range("E46").Value = "  +0.95%  "

# Row 47
$ws.Range("E47").Value = "  +2.42%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5096"
$ws.Range("E48").Value = "  -0.18%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4076"
$ws.Range("E49").Value = "  +0.46%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07478"
$ws.Range("E50").Value = "  +7.29%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.973"
$ws.Range("E51").Value = "  +1.04%  "
